# Presentation Edit, Re-schedule, View Result Function - DEVELOP
#
# 1) Re-schedule the presentation due date / date-time for rows 2-19 on the
#    "Presentation" sheet (columns E = presentationDueDate, F = presentationDateTime).
# 2) Update the selection on the "ProjectModule" sheet (view result range).
# 3) Update the selection on the "Presentation" sheet (view result range) and
#    leave it as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Presentation sheet: re-schedule rows 2-19 ---
$wsPresentation = $wb.Worksheets.Item("Presentation")

# Set column F (presentationDateTime) before column E (presentationDueDate)
# so shared-string ordering matches the source workbook.
$wsPresentation.Range("F2:F19").Value = "2024-01-01 00:00:00"
$wsPresentation.Range("E2:E19").Value = "2024-04-09 00:00:00"

# --- ProjectModule sheet: update the viewed/selected result range ---
$wsProjectModule = $wb.Worksheets.Item("ProjectModule")
$wsProjectModule.Activate()
$wsProjectModule.Range("G2:H12").Select()

# --- Presentation sheet: update the viewed/selected result range and ---
# --- re-activate it so it remains the workbook's active tab/sheet.   ---
$wsPresentation.Activate()
$wsPresentation.Range("H3:H29").Select()
